# Update translation evaluation results: filenames were reshuffled across rows,
# and Impact values for rows 15-78 were updated accordingly (mostly Correct -> Test Failed,
# with the "Infinite Loop" marker moving from row 15 to row 78).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "codeforces_99_A.py"
$ws.Range("D4").Value = "atcoder_ABC132_F.py"
$ws.Range("D5").Value = "codeforces_334_A.py"
$ws.Range("D6").Value = "atcoder_ARC102_C.py"
$ws.Range("D7").Value = "codeforces_651_A.py"
$ws.Range("D8").Value = "atcoder_ABC149_C.py"
$ws.Range("D9").Value = "atcoder_ABC155_E.py"
$ws.Range("D10").Value = "codeforces_581_A.py"
$ws.Range("D11").Value = "codeforces_58_B.py"
$ws.Range("D12").Value = "codeforces_171_A.py"
$ws.Range("D13").Value = "atcoder_ABC170_A.py"
$ws.Range("D14").Value = "atcoder_ABC158_A.py"
$ws.Range("D15").Value = "atcoder_ABC158_B.py"
$ws.Range("E15").Value = "Test Failed"
$ws.Range("D16").Value = "atcoder_ABC043_B.py"
$ws.Range("E16").Value = "Test Failed"
$ws.Range("D17").Value = "atcoder_ABC172_D.py"
$ws.Range("E17").Value = "Test Failed"
$ws.Range("D18").Value = "codeforces_546_A.py"
$ws.Range("E18").Value = "Test Failed"
$ws.Range("D19").Value = "atcoder_AGC046_B.py"
$ws.Range("E19").Value = "Test Failed"
$ws.Range("D20").Value = "codeforces_92_A.py"
$ws.Range("E20").Value = "Test Failed"
$ws.Range("D21").Value = "atcoder_ABC164_A.py"
$ws.Range("E21").Value = "Test Failed"
$ws.Range("D22").Value = "codeforces_190_A.py"
$ws.Range("E22").Value = "Test Failed"
$ws.Range("D23").Value = "atcoder_ABC114_C.py"
$ws.Range("E23").Value = "Test Failed"
$ws.Range("E24").Value = "Test Failed"
$ws.Range("D25").Value = "atcoder_ABC169_C.py"
$ws.Range("E25").Value = "Test Failed"
$ws.Range("D26").Value = "codeforces_672_A.py"
$ws.Range("E26").Value = "Test Failed"
$ws.Range("D27").Value = "codeforces_306_A.py"
$ws.Range("E27").Value = "Test Failed"
$ws.Range("D28").Value = "atcoder_ABC051_A.py"
$ws.Range("E28").Value = "Test Failed"
$ws.Range("D29").Value = "atcoder_ARC062_B.py"
$ws.Range("E29").Value = "Test Failed"
$ws.Range("D30").Value = "atcoder_ABC122_D.py"
$ws.Range("E30").Value = "Test Failed"
$ws.Range("D31").Value = "codeforces_622_A.py"
$ws.Range("E31").Value = "Test Failed"
$ws.Range("D32").Value = "codeforces_59_A.py"
$ws.Range("E32").Value = "Test Failed"
$ws.Range("D33").Value = "atcoder_ABC124_C.py"
$ws.Range("E33").Value = "Test Failed"
$ws.Range("D34").Value = "codeforces_369_B.py"
$ws.Range("E34").Value = "Test Failed"
$ws.Range("D35").Value = "codeforces_8_B.py"
$ws.Range("E35").Value = "Test Failed"
$ws.Range("D36").Value = "atcoder_ABC125_A.py"
$ws.Range("E36").Value = "Test Failed"
$ws.Range("D37").Value = "atcoder_ABC174_C.py"
$ws.Range("E37").Value = "Test Failed"
$ws.Range("D38").Value = "atcoder_ABC168_C.py"
$ws.Range("E38").Value = "Test Failed"
$ws.Range("D39").Value = "codeforces_55_A.py"
$ws.Range("E39").Value = "Test Failed"
$ws.Range("D40").Value = "codeforces_340_A.py"
$ws.Range("E40").Value = "Test Failed"
$ws.Range("D41").Value = "atcoder_ABC127_B.py"
$ws.Range("E41").Value = "Test Failed"
$ws.Range("D42").Value = "codeforces_514_A.py"
$ws.Range("E42").Value = "Test Failed"
$ws.Range("D43").Value = "atcoder_ABC178_A.py"
$ws.Range("E43").Value = "Test Failed"
$ws.Range("D44").Value = "atcoder_ABC153_A.py"
$ws.Range("E44").Value = "Test Failed"
$ws.Range("D45").Value = "codeforces_79_A.py"
$ws.Range("E45").Value = "Test Failed"
$ws.Range("D46").Value = "codeforces_49_A.py"
$ws.Range("E46").Value = "Test Failed"
$ws.Range("D47").Value = "codeforces_579_A.py"
$ws.Range("E47").Value = "Test Failed"
$ws.Range("D48").Value = "atcoder_ABC042_A.py"
$ws.Range("E48").Value = "Test Failed"
$ws.Range("D49").Value = "atcoder_AGC046_A.py"
$ws.Range("E49").Value = "Test Failed"
$ws.Range("D50").Value = "codeforces_242_A.py"
$ws.Range("E50").Value = "Test Failed"
$ws.Range("D51").Value = "atcoder_ABC124_A.py"
$ws.Range("E51").Value = "Test Failed"
$ws.Range("D52").Value = "codeforces_110_B.py"
$ws.Range("E52").Value = "Test Failed"
$ws.Range("D53").Value = "atcoder_ABC120_C.py"
$ws.Range("E53").Value = "Test Failed"
$ws.Range("D54").Value = "atcoder_AGC002_A.py"
$ws.Range("E54").Value = "Test Failed"
$ws.Range("D55").Value = "atcoder_ABC143_A.py"
$ws.Range("E55").Value = "Test Failed"
$ws.Range("D56").Value = "codeforces_678_A.py"
$ws.Range("E56").Value = "Test Failed"
$ws.Range("D57").Value = "atcoder_ABC136_B.py"
$ws.Range("E57").Value = "Test Failed"
$ws.Range("D58").Value = "atcoder_ABC178_B.py"
$ws.Range("E58").Value = "Test Failed"
$ws.Range("D59").Value = "codeforces_379_A.py"
$ws.Range("E59").Value = "Test Failed"
$ws.Range("D60").Value = "codeforces_276_B.py"
$ws.Range("E60").Value = "Test Failed"
$ws.Range("D61").Value = "codeforces_189_A.py"
$ws.Range("E61").Value = "Test Failed"
$ws.Range("D62").Value = "atcoder_ABC070_B.py"
$ws.Range("E62").Value = "Test Failed"
$ws.Range("D63").Value = "codeforces_373_B.py"
$ws.Range("E63").Value = "Test Failed"
$ws.Range("D64").Value = "atcoder_ABC169_D.py"
$ws.Range("E64").Value = "Test Failed"
$ws.Range("D65").Value = "codeforces_569_A.py"
$ws.Range("E65").Value = "Test Failed"
$ws.Range("D66").Value = "atcoder_ABC151_A.py"
$ws.Range("E66").Value = "Test Failed"
$ws.Range("D67").Value = "codeforces_32_B.py"
$ws.Range("E67").Value = "Test Failed"
$ws.Range("D68").Value = "codeforces_86_A.py"
$ws.Range("E68").Value = "Test Failed"
$ws.Range("D69").Value = "codeforces_678_B.py"
$ws.Range("E69").Value = "Test Failed"
$ws.Range("D70").Value = "atcoder_AGC025_A.py"
$ws.Range("E70").Value = "Test Failed"
$ws.Range("D71").Value = "atcoder_ABC139_B.py"
$ws.Range("E71").Value = "Test Failed"
$ws.Range("D72").Value = "codeforces_96_B.py"
$ws.Range("E72").Value = "Test Failed"
$ws.Range("D73").Value = "codeforces_459_A.py"
$ws.Range("E73").Value = "Test Failed"
$ws.Range("D74").Value = "atcoder_ABC142_A.py"
$ws.Range("E74").Value = "Test Failed"
$ws.Range("D75").Value = "atcoder_ABC149_B.py"
$ws.Range("E75").Value = "Test Failed"
$ws.Range("D76").Value = "codeforces_30_A.py"
$ws.Range("E76").Value = "Test Failed"
$ws.Range("D77").Value = "codeforces_544_B.py"
$ws.Range("E77").Value = "Test Failed"
$ws.Range("D78").Value = "codeforces_669_A.py"
$ws.Range("E78").Value = "Infinite Loop"
$ws.Range("D79").Value = "atcoder_AGC006_B.py"
$ws.Range("D80").Value = "codeforces_203_A.py"
$ws.Range("D81").Value = "codeforces_147_A.py"
$ws.Range("D82").Value = "atcoder_ABC132_A.py"
